$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("indicators")

# Add a new row (row 5) describing the new "density_dead_wood" indicator
$ws.Range("A5").Value = "density_dead_wood"
$ws.Range("F5").Value = "state, dbh, h, n"
$ws.Range("I5").Value = "Density of large dead wood. By default, large dead tree is as DBH >= 17.5 cm"
$ws.Range("J5").Value = "Núria"

# Update selection to the new row's first cell, matching the saved workbook state
$ws.Range("A5").Select()

$wb.Save()
